$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet 1): update Property/Value table ---
$meta = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: new publish date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# The old "Contact" / "No display for ContactDetail" row (row 10) becomes
# "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The second duplicate "Contact" row (row 11) is removed entirely, shifting
# every row below it up by one (Description ends up on row 11, etc.)
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements" (sheet 2): update the root Extension row's Short/Definition ---
$elements = $wb.Worksheets.Item(2)

$elements.Range("K2").Value = "CareGapComplianceEventDisplayName"
$elements.Range("L2").Value = "Text describing the treatment or compliance event required to close the care gap.  Should be suitable for display to patient."
